# Configuração Inicial.xlsx — start of part 3 use-case edits.
# 3.1 "Indica as alterações a efetuar" (B15) and 3.2 "Regressa ao passo 1"
# (C16) are merged/replaced by a single "3.1 Regressa a 2" note placed in C15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").ClearContents()
$ws.Range("C15").Value = "3.1 Regressa a 2"
$ws.Range("C16").ClearContents()

# Leave the cursor where the author left it while working on this area.
$ws.Range("C16").Select() | Out-Null
